# Applies updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.877.30"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.284.45"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.15%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "268.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.68"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("E9").Value = "  +2.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.48"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.97"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.18%  "

$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.615.23"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.30"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.825"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.270.62"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.872.69"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("E20").Value = "  +3.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.83"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.43%  "

$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +12.65%  "

$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.28"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.58%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("E26").Value = "  +5.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +9.38%  "

$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.32%  "

$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.35"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.38%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.03"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0903"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.02%  "

$ws.Range("E35").Value = "  +1.11%  "

$ws.Range("E36").Value = "  -1.11%  "

$ws.Range("E37").Value = "  -1.91%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0348"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.41"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +16.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.244"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +21.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.25"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.45"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.94"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.24"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +10.26%  "

$ws.Range("E46").Value = "  +4.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "99.40"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("E50").Value = "  -3.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.495.16"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.54%  "
